$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.973.16'
$ws.Range('D3').Value = '2.419.31'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '563.17'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = '142.79'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('E11').Value = '  -3.90%  '
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').Value = '25.95'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').Value = '0.0000173'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '2.855.02'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '61.830.01'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '2.416.14'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '11.34'
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('D19').Value = '323.31'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '66.65'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '8.72'
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('D26').Value = '549.97'
$ws.Range('E26').Value = '  -7.73%  '
$ws.Range('D27').Value = '2.538.10'
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = '4.74'
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '0.379'
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').Value = '153.88'
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('E39').Value = '  -5.03%  '
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('E41').Value = '  -0.79%  '
# D42 has a trailing zero ("0.990") that Excel would drop if it
# auto-detects the value as a Number, so force text format first.
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.990'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').Value = '146.64'
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('D44').Value = '2.23'
$ws.Range('E44').Value = '  -5.27%  '
$ws.Range('D45').Value = '3.64'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('E46').Value = '  -2.37%  '
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('E48').Value = '  +0.01%  '
# D49 has a trailing zero ("0.0920") that Excel would drop if it
# auto-detects the value as a Number, so force text format first.
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0920'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('E51').Value = '  +0.68%  '
